# Update "paises.xlsx" (sheet "Pais") with the new data pull taken at 17:52
# (instead of 17:22) and re-sort the country table descending by "Casos
# totales" (column B), which is how the source refreshes this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 17:52"

# 2) Updated per-country figures (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes).
#    Only these countries' numbers changed between the two pulls; everything
#    else keeps its previous values and simply gets re-sorted below.
$updates = @{
    "Estados Unidos" = @(1013508, 3152, 139422, 817037, 14187, 252, 57049)
    "Canada"         = @(49025,   525,  18268,  27991,  557,   59,  2766)
    "Singapur"       = @(14951,   528,  1128,   13809,  20,    0,   14)
    "Chile"          = @(14365,   552,  7710,   6448,   426,   9,   207)
    "Polonia"        = @(12218,   316,  2655,   8967,   160,   34,  596)
    "Argelia"        = @(3649,    132,  1651,   1561,   40,    5,   437)
    "Cuba"           = @(1437,    48,   575,    804,    12,    2,   58)
    "Mali"           = @(424,     16,   122,    278,    0,     1,   24)
    "Montenegro"     = @(321,     0,    199,    115,    7,     0,   7)
    "Sudan"          = @(318,     43,   31,     262,    0,     3,   25)
    "Liberia"        = @(141,     17,   45,     80,     0,     4,   16)
}

$dataRange = $ws.Range("A4:A216")
$xlWhole = 1

foreach ($name in $updates.Keys) {
    $found = $dataRange.Find($name, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $xlWhole)
    if ($found -ne $null) {
        $r = $found.Row
        $vals = $updates[$name]
        $ws.Cells.Item($r, 2).Value = $vals[0]
        $ws.Cells.Item($r, 3).Value = $vals[1]
        $ws.Cells.Item($r, 4).Value = $vals[2]
        $ws.Cells.Item($r, 5).Value = $vals[3]
        $ws.Cells.Item($r, 6).Value = $vals[4]
        $ws.Cells.Item($r, 7).Value = $vals[5]
        $ws.Cells.Item($r, 8).Value = $vals[6]
    }
}

# 3) Re-sort the whole table (rows 4-216) descending by "Casos totales".
$sortRange = $ws.Range("A4:H216")
$key1 = $ws.Range("B4:B216")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1, 0, 2, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# 4) The sort above is stable (ties keep their pre-sort row order). In the
#    refreshed source data, a row whose figures were just updated (Liberia)
#    sorts ahead of an unchanged row it now ties with (Gibraltar, both 141),
#    even though Gibraltar originally preceded it. Fix up any such adjacent
#    tie so updated rows lead ties against untouched rows, matching the
#    refreshed source ordering.
for ($r = 4; $r -lt 216; $r++) {
    $nameHere = $ws.Cells.Item($r, 1).Value2
    $nameNext = $ws.Cells.Item($r + 1, 1).Value2
    $bHere = $ws.Cells.Item($r, 2).Value2
    $bNext = $ws.Cells.Item($r + 1, 2).Value2
    if ($bHere -eq $bNext -and -not $updates.ContainsKey($nameHere) -and $updates.ContainsKey($nameNext)) {
        $rowHere = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 8)).Value2
        $rowNext = $ws.Range($ws.Cells.Item($r + 1, 1), $ws.Cells.Item($r + 1, 8)).Value2
        $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 8)).Value = $rowNext
        $ws.Range($ws.Cells.Item($r + 1, 1), $ws.Cells.Item($r + 1, 8)).Value = $rowHere
    }
}
